$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$ws.Range("E34").Value = 15
$ws.Range("F34").Value = 7
$ws.Range("H34").Value = 7

$ws.Range("E63").Value = 27
$ws.Range("F63").Value = 8
$ws.Range("H63").Value = 8

$ws.Range("E76").Value = 43

$ws.Range("E85").Value = 5
$ws.Range("F85").Value = 3
$ws.Range("H85").Value = 3
